$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 175, shifting rows 175:179 down to 176:180
$ws.Rows.Item(175).Insert()

# Populate the newly inserted row 175 with the new weekly data entry
$ws.Cells.Item(175, 1).Value = 1
$ws.Cells.Item(175, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(175, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(175, 4).Value = 44448
$ws.Cells.Item(175, 5).Value = 15
$ws.Cells.Item(175, 6).Value = 100114013
$ws.Cells.Item(175, 7).Value = "Zanahoria"
$ws.Cells.Item(175, 8).Value = "Sin especificar"
$ws.Cells.Item(175, 9).Value = "Primera"
$ws.Cells.Item(175, 10).Value = 100
$ws.Cells.Item(175, 11).Value = 8000
$ws.Cells.Item(175, 12).Value = 8500
$ws.Cells.Item(175, 13).Value = 8250
$ws.Cells.Item(175, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(175, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(175, 16).Value = 330
$ws.Cells.Item(175, 17).Value = 25
$ws.Cells.Item(175, 18).Value = "Hortaliza"
